$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

# ---------------------------------------------------------------------------
# 1. Update the "panel_query_time" (column F) timestamps on the "data" sheet
# ---------------------------------------------------------------------------
$data.Range("F2").Value = "2021-10-05 14:34:06.039507"
$data.Range("F3").Value = "2021-10-05 14:34:06.039515"
$data.Range("F4").Value = "2021-10-05 14:34:06.039518"
$data.Range("F5").Value = "2021-10-05 14:34:06.039521"
$data.Range("F6").Value = "2021-10-05 14:34:06.039524"
$data.Range("F7").Value = "2021-10-05 14:34:06.039527"
$data.Range("F8").Value = "2021-10-05 14:34:06.039530"

# ---------------------------------------------------------------------------
# 2. Add a new "metadata" worksheet after "data"
#    (copy the existing sheet so the page/outline properties carry over,
#    then wipe its contents and repopulate)
# ---------------------------------------------------------------------------
$data.Copy($null, $data)
$new = $wb.Worksheets.Item($wb.Worksheets.Count)
$new.Name = "metadata"
$new.Cells.Clear()

# Header row
$new.Range("B1").Value = "data_name"
$new.Range("C1").Value = "data_id"
$new.Range("D1").Value = "data_version"
$new.Range("E1").Value = "data_version_created"
$new.Range("F1").Value = "panel_query_time"
$new.Range("G1").Value = "panel_get_request"

# Reuse the same header style ("data" sheet's B1:F1 style) for the header row
$data.Range("B1:F1").Copy()
$new.Range("B1:F1").PasteSpecial(-4122)
$data.Range("B1").Copy()
$new.Range("G1").PasteSpecial(-4122)

# Reuse the same style as "data" sheet's A2 (bold/bordered numeric index cell)
$data.Range("A2").Copy()
$new.Range("A2").PasteSpecial(-4122)
$new.Range("A2").Value = 0

# Data row
$new.Range("B2").Value = "Hyperthyroidism"
$new.Range("C2").Value = 3372

# Force "0.19" to be stored as text, not a number, then drop the temporary
# number-format style so the cell ends up with no style index (like the rest
# of row 2).
$new.Range("D2").NumberFormat = "@"
$new.Range("D2").Value = "0.19"
$new.Range("D2").ClearFormats()

$new.Range("E2").Value = "2021-07-08T02:06:43.166637Z"
$new.Range("F2").Value = "2021-10-05 14:34:06.035845"
$new.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/3372/?format=json"

# Keep "data" as the active/selected sheet (unchanged from the original
# workbook view) now that "metadata" has been appended after it.
$data.Activate()

Write-Output "done"
